# Updated cryptos list on Thu Aug 10 06:52:47 UTC 2023 with GitHub Actions
#
# The "Price" (D) and "Volume(1h)" (E) columns hold plain text (many prices
# look numeric, e.g. "242.68", but some have two dots, e.g. "29.525.10", so
# the whole column is stored as text). Rows 49/50 additionally swap which
# coin (name/link) occupies that rank.
#
# Because several new price strings parse as plain numbers (e.g. "242.40"),
# assigning them straight to .Value would make Excel auto-convert them to
# numeric cells (losing the literal text, e.g. "242.40" -> 242.4). To keep
# them as literal text we briefly force the cell to Text format before the
# assignment (so the engine stores the literal string) and then clear the
# formatting again so the cell's style is left untouched, matching the
# original (unstyled) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

# Bitcoin
Set-TextValue "D2" "29.508.68"
Set-TextValue "E2" "  -0.65%  "
# Ethereum
Set-TextValue "D3" "1.849.68"
Set-TextValue "E3" "  -0.29%  "
# TetherUSD
Set-TextValue "D4" "0.9996"
Set-TextValue "E4" "  +0.05%  "
# BNB
Set-TextValue "D5" "242.40"
Set-TextValue "E5" "  -0.83%  "
# XRP
Set-TextValue "D6" "0.6289"
Set-TextValue "E6" "  -1.72%  "
# USDC
Set-TextValue "E7" "  +0.07%  "
# OKB
Set-TextValue "D8" "47.91"
Set-TextValue "E8" "  +1.29%  "
# Dogecoin
Set-TextValue "D9" "0.07538"
Set-TextValue "E9" "  +0.50%  "
# Cardano
Set-TextValue "D10" "0.2972"
Set-TextValue "E10" "  +0.05%  "
# Solana
Set-TextValue "D11" "24.25"
Set-TextValue "E11" "  -0.40%  "
# TRON
Set-TextValue "D12" "0.07677"
Set-TextValue "E12" "  +0.29%  "
# WrappedEther
Set-TextValue "D13" "1.875.92"
Set-TextValue "E13" "  +0.78%  "
# Polkadot
Set-TextValue "D14" "5.018"
Set-TextValue "E14" "  -0.26%  "
# Polygon
Set-TextValue "D15" "0.6849"
Set-TextValue "E15" "  -0.68%  "
# Litecoin
Set-TextValue "D16" "83.72"
Set-TextValue "E16" "  -0.03%  "
# ShibaInu
Set-TextValue "D17" "0.000009815"
Set-TextValue "E17" "  +0.21%  "
# WrappedliquidstakedEther2.0
Set-TextValue "D18" "2.132.19"
Set-TextValue "E18" "  +1.18%  "
# Uniswap
Set-TextValue "E19" "  +2.19%  "
# WrappedBTC
Set-TextValue "D20" "29.559.33"
Set-TextValue "E20" "  -0.51%  "
# BitcoinCash
Set-TextValue "E21" "  -0.77%  "
# Avalanche
Set-TextValue "D22" "12.50"
Set-TextValue "E22" "  -1.15%  "
# Chainlink
Set-TextValue "D24" "7.603"
Set-TextValue "E24" "  +1.96%  "
# BinanceUSD
Set-TextValue "E25" "  +0.07%  "
# Monero
Set-TextValue "D26" "155.75"
Set-TextValue "E26" "  -1.54%  "
# Stellar
Set-TextValue "D27" "0.1387"
Set-TextValue "E27" "  -1.96%  "
# Cosmos
Set-TextValue "D28" "8.422"
Set-TextValue "E28" "  -1.15%  "
# EthereumClassic
Set-TextValue "D29" "17.71"
Set-TextValue "E29" "  -0.99%  "
# PancakeSwap
Set-TextValue "D30" "1.481"
Set-TextValue "E30" "  -0.67%  "
# Hedera
Set-TextValue "D31" "0.05838"
Set-TextValue "E31" "  -6.09%  "
# Toncoin
Set-TextValue "D32" "1.285"
Set-TextValue "E32" "  +0.55%  "
# Filecoin
Set-TextValue "D33" "4.109"
Set-TextValue "E33" "  -1.15%  "
# InternetComputer(DFINITY)
Set-TextValue "D34" "4.040"
Set-TextValue "E34" "  -1.27%  "
# LidoDAOToken
Set-TextValue "D35" "1.895"
Set-TextValue "E35" "  +0.13%  "
# ARBITRUM
Set-TextValue "D36" "1.170"
Set-TextValue "E36" "  -0.15%  "
# ImmutableX
Set-TextValue "D37" "0.7154"
Set-TextValue "E37" "  -1.68%  "
# HuobiToken
Set-TextValue "D38" "2.590"
Set-TextValue "E38" "  -0.65%  "
# MXToken
Set-TextValue "D39" "2.801"
Set-TextValue "E39" "  -0.95%  "
# Maker
Set-TextValue "D40" "1.237.17"
Set-TextValue "E40" "  +3.02%  "
# VeChain
Set-TextValue "D41" "0.01778"
Set-TextValue "E41" "  -0.15%  "
# TrustWalletToken
Set-TextValue "D42" "0.9112"
Set-TextValue "E42" "  -1.12%  "
# FraxShare
Set-TextValue "D43" "6.134"
Set-TextValue "E43" "  -1.00%  "
# RocketPoolETH
Set-TextValue "D44" "2.041.16"
Set-TextValue "E44" "  +1.26%  "
# PaxDollar
Set-TextValue "E45" "  +0.03%  "
# Quant
Set-TextValue "D46" "101.96"
Set-TextValue "E46" "  +0.04%  "
# Aave
Set-TextValue "D47" "67.45"
Set-TextValue "E47" "  +1.54%  "
# Aptos
Set-TextValue "D48" "7.267"
Set-TextValue "E48" "  +8.92%  "

# Row 49 now holds EnergySwap (was BabyDogeCoin)
Set-TextValue "B49" "EnergySwap"
Set-TextValue "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "9.163"
Set-TextValue "E49" "  +0.10%  "

# Row 50 now holds BabyDogeCoin (was EnergySwap)
Set-TextValue "B50" "BabyDogeCoin"
Set-TextValue "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D50" "0.00000000117"
Set-TextValue "E50" "  -1.22%  "

# TheSandbox
Set-TextValue "D51" "0.4028"
Set-TextValue "E51" "  -0.69%  "

Write-Host "Applied cryptos update"
